$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The daily price-history table lives in A1:D15 (header in row 1, newest
# date on top). A new day's reading (2025-12-05) needs to be added at the
# top, pushing every existing row down by one and growing the table to
# A1:D16 - same price figures as before (783.5 / 1112 / 3610), only the
# date column changes.

$newDate = "2025-12-05"
$ironOre = 783.5
$cokingCoal = 1112
$hBeam = 3610

# Push existing data rows (2..15) down one row, keeping their formatting.
$ws.Rows.Item(2).Insert()

# Inserting a row copies the formatting of the row above (the bold/bordered
# header), so strip that back off - the data rows in this sheet are
# unstyled.
$ws.Rows.Item(2).ClearFormats()

$dateCell = $ws.Cells.Item(2, 1)
# Force text storage so the date stays a literal string like "2025-12-05"
# instead of being auto-parsed into a date serial number, then drop the
# number-format residue so the cell ends up unstyled, matching its peers.
$dateCell.NumberFormat = "@"
$dateCell.Value = $newDate
$dateCell.ClearFormats()

$ws.Cells.Item(2, 2).Value = $ironOre
$ws.Cells.Item(2, 3).Value = $cokingCoal
$ws.Cells.Item(2, 4).Value = $hBeam
